$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Phút hành chính" column (column K) entirely, shifting
# subsequent columns left - matches the commit "bỏ cột Phút HC các bảng
# công Chủ nhật" (remove the Admin-Minutes column).
$ws.Columns.Item(11).Delete()

# Reflect the resulting selection/view state: the just-deleted column's
# header cell becomes the active cell with the whole column selected,
# and the view is scrolled so column E is the left-most visible column.
$ws.Application.ActiveWindow.ScrollColumn = 5
$ws.Range("K1:K1048576").Select()
